# Updates the cryptos price/volume table (sheet1, rows 2-51) to the latest
# scraped values. Cell values are plain/formatted numeric-looking text
# (e.g. "219.35", "88.447.14", "  +9.30%  ") stored as strings in the
# source workbook, not real numbers, so values that Excel would otherwise
# auto-coerce to a number are written with a leading apostrophe to force
# text storage, matching the original "General" number format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '89.084.82'
$ws.Cells.Item(2, 5).Value = '  +9.82%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.346.90'
$ws.Cells.Item(3, 5).Value = '  +6.50%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.13%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''219.35'
$ws.Cells.Item(5, 5).Value = '  +5.87%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''647.91'
$ws.Cells.Item(6, 5).Value = '  +4.64%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''0.403'
$ws.Cells.Item(7, 5).Value = '  +44.58%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.998'
$ws.Cells.Item(8, 5).Value = '  -0.10%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.605'
$ws.Cells.Item(9, 5).Value = '  +4.87%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '3.337.27'
$ws.Cells.Item(10, 5).Value = '  +6.18%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.594'
$ws.Cells.Item(11, 5).Value = '  +3.66%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''0.0000289'
$ws.Cells.Item(12, 5).Value = '  +14.79%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''36.56'
$ws.Cells.Item(13, 5).Value = '  +17.42%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +2.42%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '3.955.05'
$ws.Cells.Item(15, 5).Value = '  +6.25%  '

# Row 16
$ws.Cells.Item(16, 5).Value = '  +5.13%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '88.817.37'
$ws.Cells.Item(17, 5).Value = '  +9.79%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '3.331.48'
$ws.Cells.Item(18, 5).Value = '  +5.98%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +6.93%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''3.14'
$ws.Cells.Item(20, 5).Value = '  -0.19%  '

# Row 21
$ws.Cells.Item(21, 2).Value = 'BitcoinCash'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Cells.Item(21, 4).Value = '''457.82'
$ws.Cells.Item(21, 5).Value = '  +7.15%  '

# Row 22
$ws.Cells.Item(22, 2).Value = 'Uniswap'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(22, 4).Value = '''9.54'
$ws.Cells.Item(22, 5).Value = '  +7.21%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''5.50'
$ws.Cells.Item(23, 5).Value = '  +9.12%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''7.47'
$ws.Cells.Item(24, 5).Value = '  +4.11%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'NEARProtocol'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(25, 4).Value = '''5.54'
$ws.Cells.Item(25, 5).Value = '  +8.35%  '

# Row 26
$ws.Cells.Item(26, 2).Value = 'Aptos'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(26, 4).Value = '''12.91'
$ws.Cells.Item(26, 5).Value = '  +20.37%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '3.480.73'
$ws.Cells.Item(27, 5).Value = '  +5.23%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'PEPE'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(28, 4).Value = '''0.0000141'
$ws.Cells.Item(28, 5).Value = '  +17.41%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Litecoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(29, 4).Value = '''78.61'
$ws.Cells.Item(29, 5).Value = '  +4.48%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''0.207'
$ws.Cells.Item(30, 5).Value = '  +42.67%  '

# Row 31
$ws.Cells.Item(31, 4).Value = '''1.00'
$ws.Cells.Item(31, 5).Value = '  -0.08%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''9.30'
$ws.Cells.Item(32, 5).Value = '  +4.82%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''595.40'
$ws.Cells.Item(33, 5).Value = '  +6.86%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''1.59'
$ws.Cells.Item(34, 5).Value = '  +8.84%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '''0.998'
$ws.Cells.Item(35, 5).Value = '  -0.28%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''2.13'
$ws.Cells.Item(36, 5).Value = '  +8.29%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''7.37'
$ws.Cells.Item(37, 5).Value = '  +25.30%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -4.75%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''23.35'
$ws.Cells.Item(39, 5).Value = '  +3.76%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''2.18'
$ws.Cells.Item(40, 5).Value = '  +9.58%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''0.427'
$ws.Cells.Item(41, 5).Value = '  +5.65%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''21.85'
$ws.Cells.Item(42, 5).Value = '  +5.47%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''0.996'
$ws.Cells.Item(43, 5).Value = '  -0.24%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''3.18'
$ws.Cells.Item(44, 5).Value = '  +6.48%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +11.64%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.05%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''157.43'
$ws.Cells.Item(47, 5).Value = '  -1.21%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''188.32'
$ws.Cells.Item(48, 5).Value = '  +1.02%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''46.18'
$ws.Cells.Item(49, 5).Value = '  +2.28%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''4.51'
$ws.Cells.Item(50, 5).Value = '  +8.38%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +7.95%  '
